$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data entered for the 2018 row (row 2): columns H, I, J and Q, R, S ---
# These cells were previously empty; fill in the newly-reported figures.
# Re-applying "General" number format keeps them on the workbook's default
# style (no explicit style override), matching how the rest of that row's
# freshly-entered figures are stored.
$ws.Range("H2").Value = 10
$ws.Range("H2").NumberFormat = "General"

$ws.Range("I2").Value = 1
$ws.Range("I2").NumberFormat = "General"

$ws.Range("J2").Value = 0
$ws.Range("J2").NumberFormat = "General"

$ws.Range("Q2").Value = 5
$ws.Range("Q2").NumberFormat = "General"

$ws.Range("R2").Value = 1
$ws.Range("R2").NumberFormat = "General"

$ws.Range("S2").Value = 0
$ws.Range("S2").NumberFormat = "General"

# --- View / window state ---
# Gridlines back on, and move the cursor to the cell that was last worked on.
$excel.ActiveWindow.DisplayGridlines = $true
$ws.Range("I3").Select()
